# "Beggining a new module Data Structures"
#
# Locate the end of the "Pointers / Pass by Value" section (the last
# section in the document) and append a new module: a page break
# followed by a "DATA STRUCTURES" Heading 1 and an "ARRAYS" Heading 2.

$d = $word.ActiveDocument

# Find the "Pass by Value" heading; the (empty) paragraph right after it
# is the anchor we insert the new module after - it is left untouched,
# exactly like in the diff.
$findRange = $d.Content
$findRange.Find.Execute("Pass by Value", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $findRange.Paragraphs(1)
$anchor = $headingPara.Next()

# Pre-create two clean paragraphs after the anchor (they inherit the
# anchor's plain "<w:rPr><w:lang .../></w:rPr>" pPr, with no direct
# spacing/indent formatting) - these will become the two heading
# paragraphs (C and D).
$anchor.Range.InsertParagraphAfter()
$afterAnchor = $anchor.Next()
$afterAnchor.Range.InsertParagraphAfter()

# Turn the anchor's trailing mark into its own paragraph carrying a
# manual page break (B). Calling InsertBreak on the still-untouched
# anchor paragraph splits it in two: the anchor itself stays empty and
# untouched, and a new paragraph holding the break is inserted right
# after it (inheriting the anchor's plain formatting, which we then
# adjust to match the target spacing/indent).
$anchor.Range.InsertBreak(7)
$breakPara = $anchor.Next()
$breakPara.LineSpacingRule = 5
$breakPara.LineSpacing = 18
$breakPara.FirstLineIndent = 0
$breakPara.Range.LanguageID = "en-US"

# "DATA STRUCTURES" - Heading 1
$h1 = $breakPara.Next()
$h1.Range.Text = "DATA STRUCTURES"
$h1.Style = "Heading 1"

# "ARRAYS" - Heading 2
$h2 = $h1.Next()
$h2.Range.Text = "ARRAYS"
$h2.Style = "Heading 2"
